$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("constants_evaluated")
$ws2 = $wb.Worksheets.Item("correlation_matrix")
$ws3 = $wb.Worksheets.Item("adj_r_squared")
$ws4 = $wb.Worksheets.Item("mol_ext_coefficients_calc")
$ws5 = $wb.Worksheets.Item("equilibrium_concentrations")
$ws6 = $wb.Worksheets.Item("absorbance_calc_abs_errors")
$ws7 = $wb.Worksheets.Item("absorbance_calc_rel_errors")

$ws1.Range("B4,C4,B5,C5").NumberFormat = "@"
$ws1.Range("B4").Value = "4.57094693183899"
$ws1.Range("C4").Value = "0.108730034369404"
$ws1.Range("B5").Value = "6.47013266563416"
$ws1.Range("C5").Value = "0.260060079433111"

$ws2.Range("B2").Value = "0.848969562068136"
$ws2.Range("A3").Value = "0.848969562068136"

$ws3.Range("A2").Value = "0.999944167785635"

$ws4.Range("C2").Value = "3.63384902787114"
$ws4.Range("D2").Value = "2964.50958508043"
$ws4.Range("E2").Value = "4768.81953495453"
$ws4.Range("F2").Value = "5054.57960446455"
$ws4.Range("G2").Value = "422942.732748066"
$ws4.Range("C3").Value = "22.5113678877864"
$ws4.Range("D3").Value = "8048.2325688426"
$ws4.Range("E3").Value = "6303.25890524787"
$ws4.Range("F3").Value = "3316.74562567488"
$ws4.Range("G3").Value = "-639021.551289459"
$ws4.Range("C4").Value = "2.25385330039595"
$ws4.Range("D4").Value = "70.2021545853265"
$ws4.Range("E4").Value = "12.6708023024833"
$ws4.Range("F4").Value = "163.831291326674"
$ws4.Range("G4").Value = "109243.711107637"
$ws4.Range("C5").Value = "1.8027054123765"
$ws4.Range("D5").Value = "56.1499739176578"
$ws4.Range("E5").Value = "10.1345211269192"
$ws4.Range("F5").Value = "131.037612580793"
$ws4.Range("G5").Value = "87376.6847412983"

$ws5.Range("A2").Value = "6.57167626841221e-08"
$ws5.Range("B2").Value = "0.000353335405410957"
$ws5.Range("C2").Value = "8.64590084348246e-07"
$ws5.Range("D2").Value = "4.50476505622661e-12"
$ws5.Range("E2").Value = "2.30315856505498e-07"
$ws5.Range("A3").Value = "8.4001348557892e-05"
$ws5.Range("B3").Value = "8.57882020990969e-05"
$ws5.Range("C3").Value = "0.000268324764177118"
$ws5.Range("D3").Value = "1.78703372421423e-06"
$ws5.Range("E3").Value = "1.8018297020471e-10"
$ws5.Range("A4").Value = "0.000347574362442624"
$ws5.Range("B4").Value = "2.46740242897292e-05"
$ws5.Range("C4").Value = "0.000319326274765001"
$ws5.Range("D4").Value = "8.7997052027367e-06"
$ws5.Range("E4").Value = "4.35464007701678e-11"
$ws5.Range("A5").Value = "0.000665648939535394"
$ws5.Range("B5").Value = "1.31508857689604e-05"
$ws5.Range("C5").Value = "0.00032594714525971"
$ws5.Range("D5").Value = "1.72019689717009e-05"
$ws5.Range("E5").Value = "2.27381305450984e-11"
$ws5.Range("A6").Value = "0.0013330579651414"
$ws5.Range("B6").Value = "6.36878430281812e-06"
$ws5.Range("C6").Value = "0.000316120388153839"
$ws5.Range("D6").Value = "3.34108307598665e-05"
$ws5.Range("E6").Value = "1.13540542723186e-11"
$ws5.Range("A7").Value = "0.00264704835806024"
$ws5.Range("B7").Value = "2.94394324161017e-06"
$ws5.Range("C7").Value = "0.000290160465876497"
$ws5.Range("D7").Value = "6.08955909026936e-05"
$ws5.Range("E7").Value = "5.71792065614301e-12"
$ws5.Range("A8").Value = "0.00664385387229889"
$ws5.Range("B8").Value = "9.35329944465505e-07"
$ws5.Range("C8").Value = "0.000231383210177002"
$ws5.Range("D8").Value = "0.000121881459947133"
$ws5.Range("E8").Value = "2.27813747491784e-12"
$ws5.Range("A9").Value = "0.0134042877665934"
$ws5.Range("B9").Value = "3.43417561645228e-07"
$ws5.Range("C9").Value = "0.000171400930360845"
$ws5.Range("D9").Value = "0.000182155652118493"
$ws5.Range("E9").Value = "1.1291620075543e-12"

$ws6.Range("C2").Value = "1.14899994949316"
$ws6.Range("D2").Value = "1.54332648404101"
$ws6.Range("E2").Value = "1.64171601946887"
$ws6.Range("F2").Value = "1.68274624679144"
$ws6.Range("G2").Value = "1.70012804176895"
$ws6.Range("H2").Value = "1.70987325013526"
$ws6.Range("I2").Value = "1.74640083440616"
$ws6.Range("J2").Value = "1.78782804934424"
$ws6.Range("C3").Value = "2.7019999511069"
$ws6.Range("D3").Value = "2.3894668422014"
$ws6.Range("E3").Value = "2.24836140197627"
$ws6.Range("F3").Value = "2.23239532649969"
$ws6.Range("G3").Value = "2.18466303856759"
$ws6.Range("H3").Value = "2.11421029064477"
$ws6.Range("I3").Value = "2.01980661505828"
$ws6.Range("J3").Value = "1.98906043905376"
$ws6.Range("C4").Value = "-5.05068364908823e-08"
$ws6.Range("D4").Value = "0.000326484041010655"
$ws6.Range("E4").Value = "-0.000283980531124994"
$ws6.Range("F4").Value = "-0.00225375320855803"
$ws6.Range("G4").Value = "-0.000871958231049019"
$ws6.Range("H4").Value = "0.00587325013525608"
$ws6.Range("I4").Value = "-0.00359916559384343"
$ws6.Range("J4").Value = "0.000828049344244297"
$ws6.Range("C5").Value = "-4.88930962383449e-08"
$ws6.Range("D5").Value = "0.000466842201398698"
$ws6.Range("E5").Value = "-0.00363859802372657"
$ws6.Range("F5").Value = "0.00439532649968877"
$ws6.Range("G5").Value = "-0.00133696143241169"
$ws6.Range("H5").Value = "0.000210290644774069"
$ws6.Range("I5").Value = "-0.000193384941721941"
$ws6.Range("J5").Value = "6.04390537617405e-05"

$ws7.Range("C2").Value = "1.14899994949316"
$ws7.Range("D2").Value = "1.54332648404101"
$ws7.Range("E2").Value = "1.64171601946887"
$ws7.Range("F2").Value = "1.68274624679144"
$ws7.Range("G2").Value = "1.70012804176895"
$ws7.Range("H2").Value = "1.70987325013526"
$ws7.Range("I2").Value = "1.74640083440616"
$ws7.Range("J2").Value = "1.78782804934424"
$ws7.Range("C3").Value = "2.7019999511069"
$ws7.Range("D3").Value = "2.3894668422014"
$ws7.Range("E3").Value = "2.24836140197627"
$ws7.Range("F3").Value = "2.23239532649969"
$ws7.Range("G3").Value = "2.18466303856759"
$ws7.Range("H3").Value = "2.11421029064477"
$ws7.Range("I3").Value = "2.01980661505828"
$ws7.Range("J3").Value = "1.98906043905376"
$ws7.Range("C4").Value = "-4.39572119154763e-08"
$ws7.Range("D4").Value = "0.000211590434874047"
$ws7.Range("E4").Value = "-0.000172947948309984"
$ws7.Range("F4").Value = "-0.00133753899617687"
$ws7.Range("G4").Value = "-0.000512615068223997"
$ws7.Range("H4").Value = "0.00344674303712211"
$ws7.Range("I4").Value = "-0.00205666605362482"
$ws7.Range("J4").Value = "0.000463374003494291"
$ws7.Range("C5").Value = "-1.80951503472779e-08"
$ws7.Range("D5").Value = "0.000195413227877228"
$ws7.Range("E5").Value = "-0.00161571848300469"
$ws7.Range("F5").Value = "0.00197276772876516"
$ws7.Range("G5").Value = "-0.000611601753161797"
$ws7.Range("H5").Value = "9.94752340463904e-05"
$ws7.Range("I5").Value = "-9.57351196643274e-05"
$ws7.Range("J5").Value = "3.03866534749827e-05"

